# Update the descr/AlternativeText metadata on the generated chart
# pictures to use the new centralized-css-utils filenames.
#
# Mapping of old -> new descr values (old descr stays unique enough to
# match on, independent of slide index):
#   slide_bc8113d8_create_sales_chart.png   -> slide_b2380a51_create_sales_chart.png
#   slide_af09c7df_create_market_share.png  -> slide_232aaf6b_create_market_share.png
#   slide_2781d4f0_create_growth_trend.png  -> slide_00cb9dae_create_growth_trend.png

$p = $ppt.ActivePresentation

$descrMap = @{
    "slide_bc8113d8_create_sales_chart.png"  = "slide_b2380a51_create_sales_chart.png";
    "slide_af09c7df_create_market_share.png" = "slide_232aaf6b_create_market_share.png";
    "slide_2781d4f0_create_growth_trend.png" = "slide_00cb9dae_create_growth_trend.png";
}

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        $oldDescr = $sh.AlternativeText
        if ($descrMap.ContainsKey($oldDescr)) {
            $sh.AlternativeText = $descrMap[$oldDescr]
        }
    }
}
